# Update the "asinh" results row (row 31) on Sheet1 with the new
# measured values, per the commit "Improved accuracy of asinh function."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 8572
$ws.Range("G31").Value = 7993
$ws.Range("M31").Value = 4500
$ws.Range("N31").Value = 3138
$ws.Range("W31").Value = 4344
$ws.Range("X31").Value = 3540

$excel.Calculate()

$wb.Save()
